$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'330.27"
$ws.Range("E2").Value = "'0.51%"
$ws.Range("G2").Value = "'23"

$ws.Range("D3").Value = "'41.15"
$ws.Range("E3").Value = "'2.52%"
$ws.Range("G3").Value = "'23"

$ws.Range("D4").Value = "'5.719"
$ws.Range("E4").Value = "'-2.19%"
$ws.Range("G4").Value = "'23"

$ws.Range("D5").Value = "'0.08128"
$ws.Range("E5").Value = "'1.25%"
$ws.Range("G5").Value = "'23"

$ws.Range("D6").Value = "'2.039"
$ws.Range("E6").Value = "'5.35%"
$ws.Range("G6").Value = "'23"

$ws.Range("D7").Value = "'8.744"
$ws.Range("E7").Value = "'0.28%"
$ws.Range("G7").Value = "'23"

$ws.Range("D8").Value = "'4.507"
$ws.Range("E8").Value = "'-1.72%"
$ws.Range("G8").Value = "'23"

$ws.Range("D9").Value = "'2.964"
$ws.Range("E9").Value = "'0.78%"
$ws.Range("G9").Value = "'23"

$ws.Range("D10").Value = "'0.9229"
$ws.Range("E10").Value = "'-2.29%"
$ws.Range("G10").Value = "'23"

$ws.Range("D11").Value = "'0.1245"
$ws.Range("E11").Value = "'-1.02%"
$ws.Range("G11").Value = "'23"

$ws.Range("D12").Value = "'0.1949"
$ws.Range("E12").Value = "'-0.88%"
$ws.Range("G12").Value = "'23"

$ws.Range("D13").Value = "'8.335"
$ws.Range("E13").Value = "'-6.28%"
$ws.Range("G13").Value = "'23"

$ws.Range("D14").Value = "'0.09326"
$ws.Range("E14").Value = "'1.67%"
$ws.Range("G14").Value = "'23"

$ws.Range("D15").Value = "'0.03656"
$ws.Range("E15").Value = "'2.95%"
$ws.Range("G15").Value = "'23"

$ws.Range("D16").Value = "'0.1054"
$ws.Range("E16").Value = "'9.41%"
$ws.Range("G16").Value = "'23"

$ws.Range("D17").Value = "'0.001307"
$ws.Range("E17").Value = "'0.15%"
$ws.Range("G17").Value = "'23"

$ws.Range("D18").Value = "'0.006189"
$ws.Range("E18").Value = "'1.42%"
$ws.Range("G18").Value = "'23"

$ws.Range("E19").Value = "'0.47%"
$ws.Range("G19").Value = "'23"

$ws.Range("D20").Value = "'0.3484"
$ws.Range("E20").Value = "'-1.22%"
$ws.Range("G20").Value = "'23"

$ws.Range("D21").Value = "'0.1415"
$ws.Range("E21").Value = "'-1.28%"
$ws.Range("G21").Value = "'23"

$ws.Range("D22").Value = "'0.2648"
$ws.Range("E22").Value = "'9.78%"
$ws.Range("G22").Value = "'23"

$ws.Range("D23").Value = "'0.04436"
$ws.Range("E23").Value = "'0.64%"
$ws.Range("G23").Value = "'23"

$ws.Range("E24").Value = "'0.13%"
$ws.Range("G24").Value = "'23"

$ws.Range("D25").Value = "'0.004298"
$ws.Range("E25").Value = "'-1.57%"
$ws.Range("G25").Value = "'23"

$ws.Range("E26").Value = "'8.70%"
$ws.Range("G26").Value = "'23"

$ws.Range("G27").Value = "'23"

$ws.Range("G28").Value = "'23"

$ws.Range("G29").Value = "'23"

$ws.Range("G30").Value = "'23"

$ws.Range("G31").Value = "'23"

$ws.Range("G32").Value = "'23"

$ws.Range("G33").Value = "'23"

$ws.Range("G34").Value = "'23"

$ws.Range("G35").Value = "'23"

$ws.Range("G36").Value = "'23"

$ws.Range("G37").Value = "'23"

$ws.Range("G38").Value = "'23"

$ws.Range("D39").Value = "'0.02825"
$ws.Range("E39").Value = "'16.28%"
$ws.Range("G39").Value = "'23"

$ws.Range("D40").Value = "'0.05496"
$ws.Range("E40").Value = "'3.95%"
$ws.Range("G40").Value = "'23"

$ws.Range("D41").Value = "'0.007637"
$ws.Range("E41").Value = "'2.47%"
$ws.Range("G41").Value = "'23"

$ws.Range("D42").Value = "'0.009940"
$ws.Range("E42").Value = "'13.96%"
$ws.Range("G42").Value = "'23"

$ws.Range("D43").Value = "'0.1424"
$ws.Range("E43").Value = "'0.36%"
$ws.Range("G43").Value = "'23"

$ws.Range("D44").Value = "'0.002121"
$ws.Range("E44").Value = "'0.97%"
$ws.Range("G44").Value = "'23"

$ws.Range("D45").Value = "'0.01175"
$ws.Range("E45").Value = "'11.45%"
$ws.Range("G45").Value = "'23"

$ws.Range("D46").Value = "'0.00006735"
$ws.Range("E46").Value = "'-2.26%"
$ws.Range("G46").Value = "'23"

$ws.Range("D47").Value = "'0.00000000749"
$ws.Range("E47").Value = "'-0.25%"
$ws.Range("G47").Value = "'23"

$ws.Range("D48").Value = "'0.002948"
$ws.Range("E48").Value = "'-6.41%"
$ws.Range("G48").Value = "'23"

$ws.Range("D49").Value = "'0.002277"
$ws.Range("E49").Value = "'60.09%"
$ws.Range("G49").Value = "'23"

$ws.Range("D50").Value = "'0.00002098"
$ws.Range("E50").Value = "'-0.25%"
$ws.Range("G50").Value = "'23"

$ws.Range("D51").Value = "'0.0001998"
$ws.Range("E51").Value = "'-0.25%"
$ws.Range("G51").Value = "'23"
